$wb = $excel.ActiveWorkbook

# --- Sheet "3.uai" was the active sheet to start with; fill in the newly
#     finished "BN_EM_POD" (row 3) results for homework5 run #3, then move
#     down to D4 the same way pressing Enter after the last edited cell would. ---
$ws3 = $wb.Worksheets.Item("3.uai")
$ws3.Activate()
$ws3.Range("B3").Value = 15844404.535230299
$ws3.Range("C3").Value = 15844404.535230299
$ws3.Range("D3").Value = 15844404.535230299

# The second chart on this sheet was a duplicate of the "BN_MLE_FOD" (row 2)
# chart; repoint it at the row that now has data (row 3, "BN_EM_POD").
$chartObj3 = $ws3.ChartObjects(2)
$series3 = $chartObj3.Chart.SeriesCollection(1)
$series3.Formula = "=SERIES('3.uai'!`$A`$3,'3.uai'!`$B`$1:`$F`$1,'3.uai'!`$B`$3:`$F`$3,1)"

$ws3.Range("D4").Select()

# --- Sheet "2.uai": same treatment. ---
$ws2 = $wb.Worksheets.Item("2.uai")
$ws2.Activate()
$ws2.Range("C3").Value = 147675298.55651101
$ws2.Range("D3").Value = 116308480.221587

$chartObj2 = $ws2.ChartObjects(2)
$series2 = $chartObj2.Chart.SeriesCollection(1)
$series2.Formula = "=SERIES('2.uai'!`$A`$3,'2.uai'!`$B`$1:`$F`$1,'2.uai'!`$B`$3:`$F`$3,1)"

$ws2.Range("D4").Select()

# --- Sheet "1.uai": finish last, leave it as the active/selected sheet. ---
$ws1 = $wb.Worksheets.Item("1.uai")
$ws1.Activate()
$ws1.Range("E3").Value = 280140818.98002499
$ws1.Range("E4").Select()
